# Household script update:
#  - rename the "DSD" sheet to "DSD_LIVESTOCK"
#  - replace the repeated "CL_COM_YESNO" codelist reference used for the
#    four livestock indicator rows with the correct per-species codelists
#  - make the DSD_LIVESTOCK sheet the active sheet/tab with F8 selected

$wb = $excel.ActiveWorkbook

$wsDsd = $wb.Worksheets.Item(1)
$wsDsd.Name = "DSD_LIVESTOCK"

$wsDsd.Range("F8").Value  = "CL_LIVESTOCK_PIG"
$wsDsd.Range("F9").Value  = "CL_LIVESTOCK_CHICKEN"
$wsDsd.Range("F10").Value = "CL_LIVESTOCK_DUCK"
$wsDsd.Range("F11").Value = "CL_LIVESTOCK_OTHER"

$wsDsd.Activate()
$wsDsd.Range("F8").Select()
